# Generate Report for Handoff
# Updates status from "In Translation" to "Ready for handoff" and refreshes
# the associated handoff timestamps across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("D2").Value = "2016-03-24 08:13:10"

# zh-cn sheet: ... | Status (C) | ... | Latest Handoff Datetime (E) | ...
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("E2").Value = "2016-03-24 08:13:05"

# de-de sheet: ... | Status (C) | ... | Latest Handoff Datetime (E) | ...
$dede.Range("C2").Value = $newStatus
$dede.Range("E2").Value = "2016-03-24 08:13:10"
